$wb = $excel.ActiveWorkbook

# --- Sheet: ChosenTreatment ---
$ws1 = $wb.Worksheets.Item("ChosenTreatment")
$ws1.Range("A2").Value = "{'SexualOrientation': 'Straight or heterosexual'}"
$ws1.Range("B2").Value = "{'DevType': 'Back-end developer'}"

# --- Sheet: Summary ---
$ws2 = $wb.Worksheets.Item("Summary")
$ws2.Range("A2").Value = 26

# --- Sheet: Subgroups ---
$ws3 = $wb.Worksheets.Item("Subgroups")
$ws3.Range("A2").Value = '{''Hobby'': ''1''}'
$ws3.Range("B2").Value = 28842
$ws3.Range("C2").Value = 5516.21828962618
$ws3.Range("D2").Value = -1197.137411587141
$ws3.Range("A3").Value = '{''Student'': ''1''}'
$ws3.Range("B3").Value = 29526
$ws3.Range("C3").Value = 5617.758312088003
$ws3.Range("D3").Value = -1095.597389125319
$ws3.Range("A4").Value = '{''UndergradMajor'': ''2''}'
$ws3.Range("B4").Value = 22173
$ws3.Range("C4").Value = 3542.149042873403
$ws3.Range("D4").Value = -3171.206658339918
$ws3.Range("A5").Value = '{''Gender'': ''1''}'
$ws3.Range("B5").Value = 33253
$ws3.Range("C5").Value = 5920.984439082945
$ws3.Range("D5").Value = -792.3712621303766
$ws3.Range("A6").Value = '{''RaceEthnicity'': ''1''}'
$ws3.Range("B6").Value = 27379
$ws3.Range("C6").Value = 6509.936391650183
$ws3.Range("D6").Value = -203.4193095631381
$ws3.Range("A7").Value = '{''Dependents'': ''2''}'
$ws3.Range("B7").Value = 24167
$ws3.Range("C7").Value = 7072.162954162452
$ws3.Range("D7").Value = 358.8072529491301
$ws3.Range("A8").Value = '{''HDI'': ''1''}'
$ws3.Range("B8").Value = 27581
$ws3.Range("C8").Value = 8280.034330059932
$ws3.Range("D8").Value = 1566.678628846611
$ws3.Range("A9").Value = '{''Hobby'': ''1'', ''Student'': ''1''}'
$ws3.Range("B9").Value = 23894
$ws3.Range("C9").Value = 4208.446115890501
$ws3.Range("D9").Value = -2504.90958532282
$ws3.Range("A10").Value = '{''Gender'': ''1'', ''Hobby'': ''1''}'
$ws3.Range("B10").Value = 27456
$ws3.Range("C10").Value = 5084.527001737077
$ws3.Range("D10").Value = -1628.828699476245
$ws3.Range("A11").Value = '{''RaceEthnicity'': ''1'', ''Hobby'': ''1''}'
$ws3.Range("B11").Value = 22529
$ws3.Range("C11").Value = 5345.204113474854
$ws3.Range("D11").Value = -1368.151587738467
$ws3.Range("A12").Value = '{''Hobby'': ''1'', ''HDI'': ''1''}'
$ws3.Range("B12").Value = 22660
$ws3.Range("C12").Value = 7025.684796894849
$ws3.Range("D12").Value = 312.3290956815272
$ws3.Range("A13").Value = '{''Gender'': ''1'', ''Student'': ''1''}'
$ws3.Range("B13").Value = 27820
$ws3.Range("C13").Value = 4860.761884746209
$ws3.Range("D13").Value = -1852.593816467112
$ws3.Range("A14").Value = '{''RaceEthnicity'': ''1'', ''Student'': ''1''}'
$ws3.Range("B14").Value = 23264
$ws3.Range("C14").Value = 5618.412897795355
$ws3.Range("D14").Value = -1094.942803417966
$ws3.Range("A15").Value = '{''HDI'': ''1'', ''Student'': ''1''}'
$ws3.Range("B15").Value = 23773
$ws3.Range("C15").Value = 7095.77558069452
$ws3.Range("D15").Value = 382.4198794811982
$ws3.Range("A16").Value = '{''Gender'': ''1'', ''UndergradMajor'': ''2''}'
$ws3.Range("B16").Value = 21070
$ws3.Range("C16").Value = 2140.616974317039
$ws3.Range("D16").Value = -4572.738726896283
$ws3.Range("A17").Value = '{''Gender'': ''1'', ''RaceEthnicity'': ''1''}'
$ws3.Range("B17").Value = 25910
$ws3.Range("C17").Value = 6471.012774252379
$ws3.Range("D17").Value = -242.3429269609423
$ws3.Range("A18").Value = '{''Gender'': ''1'', ''Dependents'': ''2''}'
$ws3.Range("B18").Value = 22568
$ws3.Range("C18").Value = 6660.009269424866
$ws3.Range("D18").Value = -53.34643178845545
$ws3.Range("A19").Value = '{''Gender'': ''1'', ''HDI'': ''1''}'
$ws3.Range("B19").Value = 25944
$ws3.Range("C19").Value = 7219.315451383841
$ws3.Range("D19").Value = 505.9597501705193
$ws3.Range("A20").Value = '{''RaceEthnicity'': ''1'', ''HDI'': ''1''}'
$ws3.Range("B20").Value = 23146
$ws3.Range("C20").Value = 7669.000695259265
$ws3.Range("D20").Value = 955.6449940459433
$ws3.Range("A21").Value = '{''Gender'': ''1'', ''Hobby'': ''1'', ''Student'': ''1''}'
$ws3.Range("B21").Value = 22758
$ws3.Range("C21").Value = 4081.320153739394
$ws3.Range("D21").Value = -2632.035547473928
$ws3.Range("A22").Value = '{''Gender'': ''1'', ''RaceEthnicity'': ''1'', ''Hobby'': ''1''}'
$ws3.Range("B22").Value = 21544
$ws3.Range("C22").Value = 5551.102665918384
$ws3.Range("D22").Value = -1162.253035294938
$ws3.Range("A23").Value = '{''Gender'': ''1'', ''Hobby'': ''1'', ''HDI'': ''1''}'
$ws3.Range("B23").Value = 21574
$ws3.Range("C23").Value = 6412.706957617019
$ws3.Range("D23").Value = -300.6487435963027
$ws3.Range("A24").Value = '{''Gender'': ''1'', ''RaceEthnicity'': ''1'', ''Student'': ''1''}'
$ws3.Range("B24").Value = 22014
$ws3.Range("C24").Value = 5982.017555487582
$ws3.Range("D24").Value = -731.3381457257392
$ws3.Range("A25").Value = '{''Gender'': ''1'', ''HDI'': ''1'', ''Student'': ''1''}'
$ws3.Range("B25").Value = 22370
$ws3.Range("C25").Value = 6176.689060932703
$ws3.Range("D25").Value = -536.6666402806186
$ws3.Range("A26").Value = '{''RaceEthnicity'': ''1'', ''HDI'': ''1'', ''Student'': ''1''}'
$ws3.Range("B26").Value = 20128
$ws3.Range("C26").Value = 6694.06933007496
$ws3.Range("D26").Value = -19.28637113836157
$ws3.Range("A27").Value = '{''Gender'': ''1'', ''RaceEthnicity'': ''1'', ''HDI'': ''1''}'
$ws3.Range("B27").Value = 21900
$ws3.Range("C27").Value = 7566.667087338407
$ws3.Range("D27").Value = 853.3113861250858
